# New rounding for VAT tables
# Round the "VAT/Sales Tax Threshold (a)" dollar amounts in column D
# from two decimal places down to whole dollars (e.g. $50,950.98 -> $50,951).
# The values are stored as literal text (not numeric currency cells), so we
# force text entry with a leading apostrophe to avoid Excel re-interpreting
# the "$"-prefixed, comma-grouped string as a numeric currency value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newThresholds = @{
    4  = '$50,951'
    5  = '$37,457'
    6  = '$31,109'
    7  = '$23,976'
    8  = '$0'
    9  = '$76,297'
    10 = '$6,908'
    11 = '$72,612'
    12 = '$11,062'
    13 = '$103,913'
    14 = '$22,456'
    15 = '$16,711'
    16 = '$57,602'
    17 = '$14,202'
    18 = '$92,218'
    19 = '$26,132'
    20 = '$90,381'
    21 = '$100,408'
    22 = '$34,205'
    23 = '$79,774'
    24 = '$97,563'
    25 = '$33,690'
    26 = '$0'
    27 = '$1,650'
    28 = '$40,813'
    29 = '$4,917'
    30 = '$111,215'
    31 = '$16,886'
    32 = '$100,763'
    33 = '$83,257'
    34 = '$0'
    35 = '$3,297'
    36 = '$81,953'
    37 = '$0'
    38 = '$119,167'
    39 = '$0'
}

foreach ($row in $newThresholds.Keys) {
    $ws.Range("D$row").Value = "'" + $newThresholds[$row]
}
